$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bsg"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.93737866666667
$ws.Range("H2").Value = 122.812136
$ws.Range("I2").Value = 0.1310914068304752
$ws.Range("J2").Value = 0.1368503435998189
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.112632333333332
$ws.Range("N2").Value = 27.337897
$ws.Range("O2").Value = 0.9981738658344552
$ws.Range("P2").Value = 0.9981738658344552
$ws.Range("Q2").Value = 373.0472804797769
$ws.Range("R2").Value = 3357.425524317992
$ws.Range("S2").Value = 0.1308520163336528
$ws.Range("T2").Value = 0.1366004365118047

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bsg"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.93737866666667
$ws.Range("H3").Value = 122.812136
$ws.Range("I3").Value = 0.1310914068304752
$ws.Range("J3").Value = 0.1368503435998189
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01667133333333333
$ws.Range("N3").Value = 0.050014
$ws.Range("O3").Value = 0.001826134165544791
$ws.Range("P3").Value = 0.001826134165544791
$ws.Range("Q3").Value = 0.6824806855448889
$ws.Range("R3").Value = 6.142326169904001
$ws.Range("S3").Value = 0.0002393904968224626
$ws.Range("T3").Value = 0.0002499070880141732

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bsg"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 141.4996693333333
$ws.Range("H4").Value = 424.499008
$ws.Range("I4").Value = 0.4531162307677896
$ws.Range("J4").Value = 0.4730219422499276
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.112632333333332
$ws.Range("N4").Value = 27.337897
$ws.Range("O4").Value = 0.9981738658344552
$ws.Range("P4").Value = 0.9981738658344552
$ws.Range("Q4").Value = 1289.434461922908
$ws.Range("R4").Value = 11604.91015730618
$ws.Range("S4").Value = 0.4522887797378216
$ws.Range("T4").Value = 0.4721581407201327

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bsg"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 141.4996693333333
$ws.Range("H5").Value = 424.499008
$ws.Range("I5").Value = 0.4531162307677896
$ws.Range("J5").Value = 0.4730219422499276
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01667133333333333
$ws.Range("N5").Value = 0.050014
$ws.Range("O5").Value = 0.001826134165544791
$ws.Range("P5").Value = 0.001826134165544791
$ws.Range("Q5").Value = 2.358988154012445
$ws.Range("R5").Value = 21.230893386112
$ws.Range("S5").Value = 0.0008274510299679385
$ws.Range("T5").Value = 0.0008638015297949479

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Bsg"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 41.22149733333333
$ws.Range("H6").Value = 123.664492
$ws.Range("I6").Value = 0.1320012236515131
$ws.Range("J6").Value = 0.1378001293072295
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.112632333333332
$ws.Range("N6").Value = 27.337897
$ws.Range("O6").Value = 0.9981738658344552
$ws.Range("P6").Value = 0.9981738658344552
$ws.Range("Q6").Value = 375.6363494281471
$ws.Range("R6").Value = 3380.727144853324
$ws.Range("S6").Value = 0.1317601717071094
$ws.Range("T6").Value = 0.137548487783085

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Bsg"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 41.22149733333333
$ws.Range("H7").Value = 123.664492
$ws.Range("I7").Value = 0.1320012236515131
$ws.Range("J7").Value = 0.1378001293072295
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01667133333333333
$ws.Range("N7").Value = 0.050014
$ws.Range("O7").Value = 0.001826134165544791
$ws.Range("P7").Value = 0.001826134165544791
$ws.Range("Q7").Value = 0.6872173225431111
$ws.Range("R7").Value = 6.184955902888
$ws.Range("S7").Value = 0.0002410519444037473
$ws.Range("T7").Value = 0.0002516415241444218

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Bsg"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.198377
$ws.Range("H8").Value = 147.595131
$ws.Range("I8").Value = 0.1575451253784747
$ws.Range("J8").Value = 0.1644661924210021
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.112632333333332
$ws.Range("N8").Value = 27.337897
$ws.Range("O8").Value = 0.9981738658344552
$ws.Range("P8").Value = 0.9981738658344552
$ws.Range("Q8").Value = 448.3267209977229
$ws.Range("R8").Value = 4034.940488979507
$ws.Range("S8").Value = 0.157257426842406
$ws.Range("T8").Value = 0.164165855087945

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Bsg"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.198377
$ws.Range("H9").Value = 147.595131
$ws.Range("I9").Value = 0.1575451253784747
$ws.Range("J9").Value = 0.1644661924210021
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01667133333333333
$ws.Range("N9").Value = 0.050014
$ws.Range("O9").Value = 0.001826134165544791
$ws.Range("P9").Value = 0.001826134165544791
$ws.Range("Q9").Value = 0.820202542426
$ws.Range("R9").Value = 7.381822881834001
$ws.Range("S9").Value = 0.0002876985360686704
$ws.Range("T9").Value = 0.0003003373330570557

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bsg"
$ws.Range("C10").Value = "Sele"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 39.4242535
$ws.Range("H10").Value = 78.848507
$ws.Range("I10").Value = 0.1262460133717474
$ws.Range("J10").Value = 0.0878613924220219
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.112632333333332
$ws.Range("N10").Value = 27.337897
$ws.Range("O10").Value = 0.9981738658344552
$ws.Range("P10").Value = 0.9981738658344552
$ws.Range("Q10").Value = 359.2587271616298
$ws.Range("R10").Value = 2155.552362969779
$ws.Range("S10").Value = 0.1260154712134654
$ws.Range("T10").Value = 0.08770094573148771

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Bsg"
$ws.Range("C11").Value = "Sele"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 39.4242535
$ws.Range("H11").Value = 78.848507
$ws.Range("I11").Value = 0.1262460133717474
$ws.Range("J11").Value = 0.0878613924220219
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01667133333333333
$ws.Range("N11").Value = 0.050014
$ws.Range("O11").Value = 0.001826134165544791
$ws.Range("P11").Value = 0.001826134165544791
$ws.Range("Q11").Value = 0.6572548715163333
$ws.Range("R11").Value = 3.943529229098
$ws.Range("S11").Value = 0.0002305421582819725
$ws.Range("T11").Value = 0.0001604466905341924
